$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "Play Beast of Wealth Slot - Free Play Included" "Play Beast of Wealth Free - Captivating Oriental Slot Game"
Replace-Text "Exciting free spins mode with jackpots" "Captivating oriental theme with four animal characters"
Replace-Text "Attractive maximum payout of €500,000" "243 ways to win and moderate volatility for steady prizes"
Replace-Text "Appealing 96.17% return to player (RTP)" "Four jackpots and maximum payout of €500,000"
Replace-Text "Moderate volatility for balanced gameplay" "Above-average theoretical RTP of 96.17%"
Replace-Text "No interactive bonus game" "Requires patience for more valuable prizes"
Replace-Text "Theme may not suit all players" "Limited variety in symbols apart from animal characters"
Replace-Text "Read our review of Beast of Wealth by Play'N'Go. Play this slot for free and enjoy exciting features, including jackpots and an RTP of 96.17%" "Play Beast of Wealth for free and experience a captivating oriental slot game with four animal characters."

Write-Output "Done"
